$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Started" Yes/No flags that changed for several players ---
$ws.Range("C8").Value  = "Yes"
$ws.Range("C9").Value  = "No"
$ws.Range("C51").Value = "Yes"
$ws.Range("C53").Value = "No"
$ws.Range("C55").Value = "No"
$ws.Range("C56").Value = "Yes"
$ws.Range("C64").Value = "No"
$ws.Range("C66").Value = "No"
$ws.Range("C69").Value = "Yes"
$ws.Range("C71").Value = "Yes"
$ws.Range("C75").Value = "No"
$ws.Range("C78").Value = "Yes"

# --- Scroll the frozen sheet back up to the top and park the selection there ---
$ws.Range("A2").Select() | Out-Null
